$wb = $excel.ActiveWorkbook

# --- Sheet1 ("URL") ---
$ws1 = $wb.Worksheets.Item("URL")

# The old duplicate "TextBoxPage" row (row 14) is replaced by what used to be
# row 15 ("radioButtonPage" / .../radio-button), and the sheet shrinks by one
# row (A1:B15 -> A1:B14): delete row 14 so row 15 shifts up into its place.
$ws1.Rows.Item(14).Delete()

# B8 (the TextBoxPage URL cell) becomes a real hyperlink, styled like the
# other Hyperlink cell (B1).
$ws1.Hyperlinks.Add($ws1.Range("B8"), "https://demoqa.com/text-box")
$ws1.Range("B8").Style = "Hyperlink"

# Selection moves to B23 and the URL sheet becomes the active tab (taking
# over from WebTables).
$ws1.Range("B23").Select()
$ws1.Activate()

# --- Sheet3 ("WebTables") ---
$ws3 = $wb.Worksheets.Item("WebTables")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1

# Restore URL as the active/selected sheet (WebTables is no longer the
# active tab in the saved workbook).
$ws1.Activate()

$wb.Save()
